$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Name row): currencies swap columns -> B2: USD -> EUR, C2: EUR -> USD
$ws.Range("B2").Value = "EUR"
$ws.Range("C2").Value = "USD"

# Row 3 (Yieldcurve row): new curve names (CMS spread cap floor update)
$ws.Range("B3").Value = "EURIBOR3M 30112020"
$ws.Range("C3").Value = "USD CF EURIBOR3M CSA"

# Row 4 (Volatility row): new vol surface names
$ws.Range("B4").Value = "EUR Vol 14Y coterm 3perc 30112020"
$ws.Range("C4").Value = "USD Vol 14Y coterm 3perc 30112020"

# Row 6 (Discountcurve row): new discount curve name
$ws.Range("B6").Value = "EURIBOR3M 30112020"
